$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "Enfermedad"
$ws.Range("C2").Value = "EPOC"
$ws.Range("C3").Value = "Obesidad"
$ws.Range("C4").Value = "Enfermedad Renal"
$ws.Range("C5").Value = "Hipertensión"
$ws.Range("C6").Value = "Cardiopatía Isquémica"
$ws.Range("C7").Value = "Enfermedad Vascular Periférica"
$ws.Range("C8").Value = "Enfermedad Valvular Cardíaca"
$ws.Range("C9").Value = "Cáncer"
$ws.Range("C10").Value = "Neumonía"
$ws.Range("C11").Value = "ACV"
